$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 89.166664
$ws.Range("I9").Value = 67
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 67
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 102
$ws.Range("N9").Value = -538

$ws.Range("H12").Value = 3250
$ws.Range("I12").Value = 4300
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 4300
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -4130
$ws.Range("N12").Value = -440

$ws.Range("H21").Value = 4600
$ws.Range("I21").Value = 4950
$ws.Range("J21").Value = 4250
$ws.Range("K21").Value = 4950
$ws.Range("L21").Value = 4250
$ws.Range("M21").Value = -4482
$ws.Range("N21").Value = -5186

$ws.Range("H23").Value = 4600
$ws.Range("I23").Value = 4950
$ws.Range("J23").Value = 4250
$ws.Range("K23").Value = 4950
$ws.Range("L23").Value = 4250
$ws.Range("M23").Value = -4716
$ws.Range("N23").Value = -4718

$ws.Range("H29").Value = 318
$ws.Range("I29").Value = 318
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 954
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -673

$ws.Range("H38").Value = 73.833336
$ws.Range("I38").Value = 73.833336
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 221.500008
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 150.499992

$ws.Range("H58").Value = 1589.5
$ws.Range("I58").Value = 1589.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4768.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -4618.5

$ws.Range("H61").Value = 800
$ws.Range("I61").Value = 800
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2400
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2228

$ws.Range("H138").Value = 3345.7693
$ws.Range("I138").Value = 623.75
$ws.Range("J138").Value = 4555.5557
$ws.Range("K138").Value = 1871.25
$ws.Range("L138").Value = 13666.6671
$ws.Range("M138").Value = 3268.75
$ws.Range("N138").Value = -23946.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H4").Value = 202
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 202
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 202
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -434

$ws.Range("H35").Value = 1176
$ws.Range("I35").Value = 1176
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1176
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -770

$ws.Range("H45").Value = 705.4
$ws.Range("I45").Value = 705.4
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 705.4
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -328.4

$ws.Range("H74").Value = 1581
$ws.Range("I74").Value = 1677.8334
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1677.8334
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -803.8334
$ws.Range("N74").Value = -2748

$ws.Range("H77").Value = 1581
$ws.Range("I77").Value = 1677.8334
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 8389.166999999999
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -4021.166999999999
$ws.Range("N77").Value = -13736

$ws.Range("H92").Value = 54250
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 54250
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 54250
$ws.Range("N92").Value = -59242

$ws.Range("H102").Value = 1488.1428
$ws.Range("I102").Value = 1488.1428
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1488.1428
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 133.8571999999999
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 61330.75
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 61330.75
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 61330.75
$ws.Range("N103").Value = -63674.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 20.416666
$ws.Range("I7").Value = 16.571428
$ws.Range("J7").Value = 25.8
$ws.Range("K7").Value = 16.571428
$ws.Range("L7").Value = 25.8
$ws.Range("M7").Value = 96.428572
$ws.Range("N7").Value = -251.8

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H31").Value = 15473.714
$ws.Range("I31").Value = 6449.8335
$ws.Range("J31").Value = 22241.625
$ws.Range("K31").Value = 6449.8335
$ws.Range("L31").Value = 22241.625
$ws.Range("M31").Value = -6154.8335
$ws.Range("N31").Value = -22831.625

$ws.Range("H34").Value = 15473.714
$ws.Range("I34").Value = 6449.8335
$ws.Range("J34").Value = 22241.625
$ws.Range("K34").Value = 6449.8335
$ws.Range("L34").Value = 22241.625
$ws.Range("M34").Value = -6247.8335
$ws.Range("N34").Value = -22645.625

$ws.Range("H99").Value = 1003398.8
$ws.Range("I99").Value = 834831.3
$ws.Range("J99").Value = 1256250
$ws.Range("K99").Value = 834831.3
$ws.Range("L99").Value = 1256250
$ws.Range("M99").Value = -833333.3
$ws.Range("N99").Value = -1259246

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 1003398.8
$ws.Range("I126").Value = 834831.3
$ws.Range("J126").Value = 1256250
$ws.Range("K126").Value = 2504493.9
$ws.Range("L126").Value = 3768750
$ws.Range("M126").Value = -2502023.9
$ws.Range("N126").Value = -3773690

$ws.Range("H134").Value = 1198.6666
$ws.Range("I134").Value = 1148.2858
$ws.Range("J134").Value = 1375
$ws.Range("K134").Value = 3444.8574
$ws.Range("L134").Value = 4125
$ws.Range("M134").Value = -909.8574000000003
$ws.Range("N134").Value = -9195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1511.7548
$ws.Range("I4").Value = 934.7273
$ws.Range("J4").Value = 2463.85
$ws.Range("K4").Value = 2804.1819
$ws.Range("L4").Value = 7391.549999999999
$ws.Range("M4").Value = -2692.1819
$ws.Range("N4").Value = -7615.549999999999

$ws.Range("H8").Value = 499
$ws.Range("I8").Value = 499
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1497
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1358

$ws.Range("H68").Value = 2276.889
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2276.889
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6830.667
$ws.Range("N68").Value = -8452.667000000001

$ws.Range("H71").Value = 2276.889
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2276.889
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 20492.001
$ws.Range("N71").Value = -28604.001

$ws.Range("H137").Value = 5990
$ws.Range("I137").Value = 7750
$ws.Range("J137").Value = 4816.6665
$ws.Range("K137").Value = 23250
$ws.Range("L137").Value = 14449.9995
$ws.Range("M137").Value = -18150
$ws.Range("N137").Value = -24649.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 60000.8
$ws.Range("I40").Value = 55000
$ws.Range("J40").Value = 63334.668
$ws.Range("K40").Value = 55000
$ws.Range("L40").Value = 63334.668
$ws.Range("M40").Value = -54864
$ws.Range("N40").Value = -63606.668

$ws.Range("H43").Value = 605375
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 605375
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 605375
$ws.Range("N43").Value = -605761

$ws.Range("H55").Value = 402.1875
$ws.Range("I55").Value = 548.5
$ws.Range("J55").Value = 381.2857
$ws.Range("K55").Value = 548.5
$ws.Range("L55").Value = 381.2857
$ws.Range("M55").Value = -375.5
$ws.Range("N55").Value = -727.2857

$ws.Range("H99").Value = 22998.715
$ws.Range("I99").Value = 22998.8
$ws.Range("J99").Value = 22998.5
$ws.Range("K99").Value = 22998.8
$ws.Range("L99").Value = 22998.5
$ws.Range("M99").Value = -20003.8
$ws.Range("N99").Value = -28988.5

$ws.Range("H100").Value = 4499
$ws.Range("I100").Value = 4499
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4499
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3958

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 27828
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 27828
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 27828
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -28280

$ws.Range("H123").Value = 80000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 80000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

$ws.Range("H126").Value = 2999
$ws.Range("I126").Value = 2999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6527
